$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.949.72"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.143.55"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.97"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.36"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.141.08"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.03"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").Value = "3.662.09"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "3.143.08"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "62.915.69"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.15"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.67"
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.91"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.04"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.75"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  -5.67%  "
$ws.Range("E34").Value = "  -4.90%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.35"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -3.96%  "
$ws.Range("D38").Value = "0.0₃0695"
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "415.18"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.73"
$ws.Range("E41").Value = "  -6.82%  "
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "2.919.58"
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("E44").Value = "  -6.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.257"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.32"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  -8.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.51"
$ws.Range("E51").Value = "  -0.51%  "
